$wb = $excel.ActiveWorkbook

$wsCompany = $wb.Worksheets.Item("Company")
$wsRecordTypes = $wb.Worksheets.Item("CompanyRecordTypes")

# Content change: revert the "new" capital-provider test company name back to the old one
$wsCompany.Range("B2").Value = "TestCapitalProvider"

# Selection / active-sheet state: CompanyRecordTypes becomes the active tab,
# and the Company sheet's remembered selection moves to J11
$wsCompany.Range("J11").Select() | Out-Null
$wsRecordTypes.Activate() | Out-Null
$wsRecordTypes.Range("G7").Select() | Out-Null
